$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing data row, pushing ticket 238's row down to row 3.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).ClearFormats()

# New ticket (241) becomes row 2.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "241"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "EH & S - Environmental Health & Safety"
$ws.Range("C2").Value = "In Progress"
$ws.Range("D2").Value = "2017-08-15T09:16:49.403000"
$ws.Range("E2").Value = "Vincent Chov"
$ws.Range("F2").Value = "Website"
$ws.Range("G2").Value = "Vincent Chov"
$ws.Range("H2").Value = "None"
$ws.Range("I2").Value = "2017-08-15T09:17:13.210000"
$ws.Range("J2").Value = "CT - PSI Hartford Office"
$ws.Range("K2").Value = "Yes"

# The previously-existing ticket 238 row (now row 3) gets an updated "EH & S Issues" timestamp.
$ws.Range("I3").Value = "2017-08-14T16:27:57.897000"
